$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) and "全部类型" (sheet4) both contain the same event
# rows; update the "想去人数" (F column) counts for the two affected rows.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1232
    $ws.Range("F6").Value = 151
}
